# Apply the "Convertation to the sheet" edit:
# Turns the numeric day-of-week header (2..6) into weekday names,
# and the numeric period index in column A (1..12) into actual class times.
# Also moves a few class entries (Desenho Técnico, EAP, Circuitos Elétricos 2)
# to their corrected day/time slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header - weekday names
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# Column A: period index -> class start time
$ws.Range("A2").Value  = "7:00"
$ws.Range("A3").Value  = "7:50"
$ws.Range("A4").Value  = "8:40"
$ws.Range("A5").Value  = "9:30"
$ws.Range("A6").Value  = "10:40"
$ws.Range("A7").Value  = "11:30"
$ws.Range("A8").Value  = "13:00"
$ws.Range("A9").Value  = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# Class entries relocated to their corrected day/time cells
$ws.Range("C5").Value  = "-"
$ws.Range("C6").Value  = "EAP"
$ws.Range("B8").Value  = "Circuitos Elétricos 2"
$ws.Range("D8").Value  = "-"
$ws.Range("D9").Value  = "-"
$ws.Range("F9").Value  = "EAP"
$ws.Range("D10").Value = "Desenho Técnico"
$ws.Range("E12").Value = "-"
